$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.024087700799368
$ws.Range("D2").Value = 1.029272626137886
$ws.Range("E2").Value = 1.049177526682449
$ws.Range("F2").Value = 1.053656323058683
$ws.Range("I2").Value = 1.031414362012518
$ws.Range("J2").Value = 1.029264277764251
$ws.Range("K2").Value = 1.032087017729598
$ws.Range("L2").Value = 1.051935310672316
$ws.Range("M2").Value = 1.056401689340812
$ws.Range("N2").Value = 1.030725951096939

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.025030107429419
$ws.Range("D3").Value = 1.029963383108125
$ws.Range("E3").Value = 1.050267894788851
$ws.Range("F3").Value = 1.054772960050087
$ws.Range("I3").Value = 1.031568386578662
$ws.Range("J3").Value = 1.029845477737807
$ws.Range("K3").Value = 1.032586201499522
$ws.Range("L3").Value = 1.052837153939144
$ws.Range("M3").Value = 1.057330637315593
$ws.Range("N3").Value = 1.031307976441122

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.025640150629473
$ws.Range("D4").Value = 1.030410349522109
$ws.Range("E4").Value = 1.050974438251363
$ws.Range("F4").Value = 1.055496334277439
$ws.Range("I4").Value = 1.031666725142655
$ws.Range("J4").Value = 1.030221202080189
$ws.Range("K4").Value = 1.032908532547807
$ws.Range("L4").Value = 1.053421114769386
$ws.Range("M4").Value = 1.057931974398212
$ws.Range("N4").Value = 1.031684234355166

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.025896670329306
$ws.Range("D5").Value = 1.030598253111887
$ws.Range("E5").Value = 1.051271708730625
$ws.Range("F5").Value = 1.05580064039514
$ws.Range("I5").Value = 1.031707748877826
$ws.Range("J5").Value = 1.030379071869416
$ws.Range("K5").Value = 1.033043878277561
$ws.Range("L5").Value = 1.053666709213174
$ws.Range("M5").Value = 1.058184834879116
$ws.Range("N5").Value = 1.031842328337602

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.025939744470666
$ws.Range("D6").Value = 1.030629802852484
$ws.Range("E6").Value = 1.051321635792349
$ws.Range("F6").Value = 1.05585174642372
$ws.Range("I6").Value = 1.031714618294958
$ws.Range("J6").Value = 1.030405573912575
$ws.Range("K6").Value = 1.033066593882967
$ws.Range("L6").Value = 1.05370795129727
$ws.Range("M6").Value = 1.058227294676977
$ws.Range("N6").Value = 1.031868868016703

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.025643578032987
$ws.Range("D7").Value = 1.030412860304872
$ws.Range("E7").Value = 1.050978409453009
$ws.Range("F7").Value = 1.055500399646918
$ws.Range("I7").Value = 1.031667274553084
$ws.Range("J7").Value = 1.030223311877847
$ws.Range("K7").Value = 1.032910341682622
$ws.Range("L7").Value = 1.053424396032427
$ws.Range("M7").Value = 1.057935352903508
$ws.Range("N7").Value = 1.031686347148979

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.024406140801631
$ws.Range("D8").Value = 1.029506069872074
$ws.Range("E8").Value = 1.049545813916905
$ws.Range("F8").Value = 1.054033522529312
$ws.Range("I8").Value = 1.031466689552298
$ws.Range("J8").Value = 1.029460769242456
$ws.Range("K8").Value = 1.032255858326476
$ws.Range("L8").Value = 1.052240008150659
$ws.Range("M8").Value = 1.056715580538799
$ws.Range("N8").Value = 1.030922721615581

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.022227508991065
$ws.Range("D9").Value = 1.027908245709976
$ws.Range("E9").Value = 1.047029090631508
$ws.Range("F9").Value = 1.051455113950199
$ws.Range("I9").Value = 1.031103099734808
$ws.Range("J9").Value = 1.028114415689477
$ws.Range("K9").Value = 1.031097447762105
$ws.Range("L9").Value = 1.050156104924829
$ws.Range("M9").Value = 1.054568085917175
$ws.Range("N9").Value = 1.029574456086128

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.020776396038064
$ws.Range("D10").Value = 1.026843137884787
$ws.Range("E10").Value = 1.045356466884498
$ws.Range("F10").Value = 1.049740510781022
$ws.Range("I10").Value = 1.030853918959771
$ws.Range("J10").Value = 1.027215102252489
$ws.Range("K10").Value = 1.030321775072134
$ws.Range("L10").Value = 1.048768965641038
$ws.Range("M10").Value = 1.053137723925747
$ws.Range("N10").Value = 1.028673865520906

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.020148365604456
$ws.Range("D11").Value = 1.026381974648319
$ws.Range("E11").Value = 1.044633437124123
$ws.Range("F11").Value = 1.048999100962146
$ws.Range("I11").Value = 1.030744416235559
$ws.Range("J11").Value = 1.026825284202858
$ws.Range("K11").Value = 1.029985103915315
$ws.Range("L11").Value = 1.048168827420961
$ws.Range("M11").Value = 1.052518674306459
$ws.Range("N11").Value = 1.028283493884933

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.019915134373113
$ws.Range("D12").Value = 1.02621068455344
$ws.Range("E12").Value = 1.044365056263345
$ws.Range("F12").Value = 1.048723862511606
$ws.Range("I12").Value = 1.030703501178726
$ws.Range("J12").Value = 1.02668042753745
$ws.Range("K12").Value = 1.029859929851905
$ws.Range("L12").Value = 1.047945985009612
$ws.Range("M12").Value = 1.052288778037203
$ws.Range("N12").Value = 1.028138431506444

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.019965161118381
$ws.Range("D13").Value = 1.026247426549543
$ws.Range("E13").Value = 1.044422616513078
$ws.Range("F13").Value = 1.048782895113531
$ws.Range("I13").Value = 1.030712288490975
$ws.Range("J13").Value = 1.026711502529098
$ws.Range("K13").Value = 1.029886785504738
$ws.Range("L13").Value = 1.047993782037144
$ws.Range("M13").Value = 1.05233808946686
$ws.Range("N13").Value = 1.028169550628144

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.020129085671621
$ws.Range("D14").Value = 1.026367815617221
$ws.Range("E14").Value = 1.044611248911859
$ws.Range("F14").Value = 1.048976346519685
$ws.Range("I14").Value = 1.030741039095636
$ws.Range("J14").Value = 1.026813311544146
$ws.Range("K14").Value = 1.029974759419686
$ws.Range("L14").Value = 1.048150405653485
$ws.Range("M14").Value = 1.052499670062608
$ws.Range("N14").Value = 1.028271504223672

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.02023009127323
$ws.Range("D15").Value = 1.026441992185632
$ws.Range("E15").Value = 1.044727496023454
$ws.Range("F15").Value = 1.049095558751167
$ws.Range("I15").Value = 1.030758721386032
$ws.Range("J15").Value = 1.026876031389186
$ws.Range("K15").Value = 1.030028947252227
$ws.Range("L15").Value = 1.048246916662131
$ws.Range("M15").Value = 1.052599231346685
$ws.Range("N15").Value = 1.028334313138085

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02081808294591
$ws.Range("D16").Value = 1.026873744614369
$ws.Range("E16").Value = 1.045404477863449
$ws.Range("F16").Value = 1.049789737316382
$ws.Range("I16").Value = 1.030861152495288
$ws.Range("J16").Value = 1.027240964620184
$ws.Range("K16").Value = 1.030344102041946
$ws.Range("L16").Value = 1.048808805435852
$ws.Range("M16").Value = 1.053178814691607
$ws.Range("N16").Value = 1.02869976461613

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.021186998218952
$ws.Range("D17").Value = 1.027144581806666
$ws.Range("E17").Value = 1.045829459416401
$ws.Range("F17").Value = 1.050225451678725
$ws.Range("I17").Value = 1.030924975175735
$ws.Range("J17").Value = 1.027469768426606
$ws.Range("K17").Value = 1.030541576725592
$ws.Range("L17").Value = 1.049161398015096
$ws.Range("M17").Value = 1.05354245427255
$ws.Range("N17").Value = 1.028928893350201

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.021402210234192
$ws.Range("D18").Value = 1.02730255990619
$ws.Range("E18").Value = 1.046077462375623
$ws.Range("F18").Value = 1.05047969534495
$ws.Range("I18").Value = 1.03096204690279
$ws.Range("J18").Value = 1.027603186273249
$ws.Range("K18").Value = 1.030656683151734
$ws.Range("L18").Value = 1.049367107747494
$ws.Range("M18").Value = 1.053754588648554
$ws.Range("N18").Value = 1.02906250066549

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.021475597039349
$ws.Range("D19").Value = 1.027356426888859
$ws.Range("E19").Value = 1.046162045054406
$ws.Range("F19").Value = 1.050566402624933
$ws.Range("I19").Value = 1.030974661107656
$ws.Range("J19").Value = 1.027648671585095
$ws.Range("K19").Value = 1.030695918348463
$ws.Range("L19").Value = 1.049437257661906
$ws.Range("M19").Value = 1.053826925994752
$ws.Range("N19").Value = 1.029108050571697

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.021147413997063
$ws.Range("D20").Value = 1.027115523184638
$ws.Range("E20").Value = 1.04578385066866
$ws.Range("F20").Value = 1.050178693424949
$ws.Range("I20").Value = 1.030918143631855
$ws.Range("J20").Value = 1.027445224029091
$ws.Range("K20").Value = 1.030520397542258
$ws.Range("L20").Value = 1.04912356315578
$ws.Range("M20").Value = 1.053503436127367
$ws.Range("N20").Value = 1.028904314096827

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.020080812668439
$ws.Range("D21").Value = 1.026332363849865
$ws.Range("E21").Value = 1.044555696274965
$ws.Range("F21").Value = 1.048919375644629
$ws.Range("I21").Value = 1.030732579405237
$ws.Range("J21").Value = 1.026783333002467
$ws.Range("K21").Value = 1.029948856584582
$ws.Range("L21").Value = 1.048104281812511
$ws.Range("M21").Value = 1.052452087323568
$ws.Range("N21").Value = 1.028241483109025

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.019410471616924
$ws.Range("D22").Value = 1.025839998453262
$ws.Range("E22").Value = 1.04378457491221
$ws.Range("F22").Value = 1.048128484653765
$ws.Range("I22").Value = 1.030614514405505
$ws.Range("J22").Value = 1.026366823889295
$ws.Range("K22").Value = 1.029588815627544
$ws.Range("L22").Value = 1.047463856930914
$ws.Range("M22").Value = 1.051791330745493
$ws.Range("N22").Value = 1.027824382505138

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.019765805994713
$ws.Range("D23").Value = 1.026101006689712
$ws.Range("E23").Value = 1.04419325963366
$ws.Range("F23").Value = 1.048547666297523
$ws.Range("I23").Value = 1.030677234845871
$ws.Range("J23").Value = 1.026587656404216
$ws.Range("K23").Value = 1.029779745349684
$ws.Range("L23").Value = 1.047803316719697
$ws.Range("M23").Value = 1.052141584802844
$ws.Range("N23").Value = 1.02804552862756

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.021165300310393
$ws.Range("D24").Value = 1.027128653514869
$ws.Range("E24").Value = 1.045804458932601
$ws.Range("F24").Value = 1.050199821161618
$ws.Range("I24").Value = 1.030921230991302
$ws.Range("J24").Value = 1.027456314708025
$ws.Range("K24").Value = 1.030529967742185
$ws.Range("L24").Value = 1.049140658950349
$ws.Range("M24").Value = 1.053521066656242
$ws.Range("N24").Value = 1.028915420525797

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.022790509669308
$ws.Range("D25").Value = 1.028321307609226
$ws.Range("E25").Value = 1.047678810277403
$ws.Range("F25").Value = 1.05212093176579
$ws.Range("I25").Value = 1.031198294646969
$ws.Range("J25").Value = 1.028462790844085
$ws.Range("K25").Value = 1.031397527286407
$ws.Range("L25").Value = 1.050694469898671
$ws.Range("M25").Value = 1.055123037166346
$ws.Range("N25").Value = 1.029923325973415
